$wb = $excel.ActiveWorkbook

# Data rows 2-14: updated "want to go" counts (column F) and new cover image
# links (column J), scraped from bilibili on the newer crawl run.
$rowData = @(
    @{ Row = 2;  F = 1648;  J = "//i1.hdslb.com/bfs/openplatform/202312/vtGcfnyc1703060683812.jpeg" },
    @{ Row = 3;  F = 222;   J = "//i2.hdslb.com/bfs/openplatform/202311/Z7mV6VXN1701160508967.jpeg" },
    @{ Row = 4;  F = 209;   J = "//i0.hdslb.com/bfs/openplatform/202311/5AgvDWGQ1700817845950.jpeg" },
    @{ Row = 5;  F = 6532;  J = "//i1.hdslb.com/bfs/openplatform/202401/OwXCPyFi1704358608332.jpeg" },
    @{ Row = 6;  F = 383;   J = "//i1.hdslb.com/bfs/openplatform/202401/bHsHJ3f21704186294427.jpeg" },
    @{ Row = 7;  F = 268;   J = "//i1.hdslb.com/bfs/openplatform/202401/VHHzVjad1704438989848.jpeg" },
    @{ Row = 8;  F = 60;    J = "//i2.hdslb.com/bfs/openplatform/202312/oPrKUOby1703664065719.jpeg" },
    @{ Row = 9;  F = 15;    J = "//i2.hdslb.com/bfs/openplatform/202401/oWbVnOjD1704445446390.jpeg" },
    @{ Row = 10; F = 8904;  J = "//i2.hdslb.com/bfs/openplatform/202312/C3P0Encm1701659824998.jpeg" },
    @{ Row = 11; F = 2364;  J = "//i1.hdslb.com/bfs/openplatform/202401/tqrMA6qB1704787264871.jpeg" },
    @{ Row = 12; F = 271;   J = "//i0.hdslb.com/bfs/openplatform/202312/X0PZ3YhH1703822037665.jpeg" },
    @{ Row = 13; F = 6660;  J = "//i2.hdslb.com/bfs/openplatform/202310/9xMTQMlg1696736126094.png" },
    @{ Row = 14; F = 10367; J = "//i2.hdslb.com/bfs/openplatform/202312/lau3mW031702535438289.jpeg" }
)

# Sheets "展览" and "全部类型" carry the full 14-row table; "演出" and
# "本地生活" only have the header row. All four get the new "Cover" column.
foreach ($sheetName in @("展览", "演出", "本地生活", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Cells.Item(1, 10).Value = "Cover"

    if ($sheetName -eq "展览" -or $sheetName -eq "全部类型") {
        foreach ($entry in $rowData) {
            $ws.Cells.Item($entry.Row, 6).Value = $entry.F
            $ws.Cells.Item($entry.Row, 10).Value = $entry.J
        }
    }
}
